$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the guest lecture cells on row 14
$ws.Range("D14").Value = "07.04: Gjesteforelesning med Ole-Petter Moe Hansen, Tryg forsikring"
$ws.Range("B14").Value = "Gjesteforelesning"

# Reflect the new selection left by the author after editing
$ws.Range("B15").Select()
